$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Update the text of the "ID 07" question (shared string change).
$ws.Range("E14").Value = "¿Está usted satisfecho con los tiempos de espera?"

# 2) Remove the empty spacer column A so the data block (originally D:I)
#    shifts left by one column to become C:H.
$ws.Range("A1").EntireColumn.Delete()

# 3) Re-apply the (manually tweaked) column widths for the now-shifted
#    data columns C:H.
$ws.Columns("C").ColumnWidth = 10.42578125
$ws.Columns("D").ColumnWidth = 23.28515625
$ws.Columns("E").ColumnWidth = 13.28515625
$ws.Columns("F").ColumnWidth = 14.28515625
$ws.Columns("G").ColumnWidth = 15.140625
$ws.Columns("H").ColumnWidth = 35

# 4) Re-apply the row heights, which were adjusted after the columns got
#    narrower and the question text re-wrapped.
$ws.Rows("4").RowHeight = 33.75
$ws.Rows("5").RowHeight = 31.5
$ws.Rows("8").RowHeight = 44.25
$ws.Rows("9").RowHeight = 34.5
$ws.Rows("10").RowHeight = 46.5
$ws.Rows("11").RowHeight = 45
$ws.Rows("12").RowHeight = 44.25
$ws.Rows("13").RowHeight = 48.75
$ws.Rows("14").RowHeight = 48
$ws.Rows("15").RowHeight = 73.5
$ws.Rows("16").RowHeight = 73.5
$ws.Rows("17").RowHeight = 63
$ws.Rows("18").RowHeight = 63.75

# 5) Match the final on-screen selection/scroll state: whole row 11 selected,
#    viewport scrolled back to the top-left (topLeftCell reset).
$ws.Rows("11").Select()
